$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44308
$ws.Range("I2").Value = 'Primera'
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 600
$ws.Range("L2").Value = 700
$ws.Range("M2").Value = 650
$ws.Range("N2").Value = '$/paquete 6 unidades'
$ws.Range("O2").Value = 'Región de Ñuble'
$ws.Range("P2").Value = 108
$ws.Range("Q2").Value = 6

$ws.Range("D3").Value = 44308
$ws.Range("I3").Value = 'Segunda'
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 500
$ws.Range("L3").Value = 500
$ws.Range("M3").Value = 500
$ws.Range("N3").Value = '$/paquete 6 unidades'
$ws.Range("O3").Value = 'Región de Ñuble'
$ws.Range("P3").Value = 83
$ws.Range("Q3").Value = 6

$ws.Range("D4").Value = 44657
$ws.Range("I4").Value = 'Primera'
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 8000
$ws.Range("L4").Value = 9000
$ws.Range("M4").Value = 8500
$ws.Range("N4").Value = '$/paquete 36 unidades'
$ws.Range("O4").Value = 'Región Metropolitana'
$ws.Range("P4").Value = 236
$ws.Range("Q4").Value = 36

$ws.Range("D5").Value = 44631
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 220
$ws.Range("K5").Value = 6000
$ws.Range("L5").Value = 6500
$ws.Range("M5").Value = 6227
$ws.Range("N5").Value = '$/paquete 36 unidades'
$ws.Range("O5").Value = 'Región Metropolitana'
$ws.Range("P5").Value = 173
$ws.Range("Q5").Value = 36

$ws.Range("D6").Value = 44321
$ws.Range("I6").Value = 'Primera'
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 600
$ws.Range("L6").Value = 700
$ws.Range("M6").Value = 650
$ws.Range("N6").Value = '$/paquete 6 unidades'
$ws.Range("O6").Value = 'Región de Ñuble'
$ws.Range("P6").Value = 108
$ws.Range("Q6").Value = 6

$ws.Range("D7").Value = 44321
$ws.Range("I7").Value = 'Segunda'
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 500
$ws.Range("L7").Value = 500
$ws.Range("M7").Value = 500
$ws.Range("N7").Value = '$/paquete 6 unidades'
$ws.Range("O7").Value = 'Región de Ñuble'
$ws.Range("P7").Value = 83
$ws.Range("Q7").Value = 6

$ws.Range("D8").Value = 44230
$ws.Range("I8").Value = 'Primera'
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 600
$ws.Range("L8").Value = 700
$ws.Range("M8").Value = 650
$ws.Range("N8").Value = '$/paquete 6 unidades'
$ws.Range("O8").Value = 'Región de Ñuble'
$ws.Range("P8").Value = 108
$ws.Range("Q8").Value = 6

$ws.Range("D9").Value = 44230
$ws.Range("I9").Value = 'Segunda'
$ws.Range("J9").Value = 50
$ws.Range("K9").Value = 500
$ws.Range("L9").Value = 500
$ws.Range("M9").Value = 500
$ws.Range("N9").Value = '$/paquete 6 unidades'
$ws.Range("O9").Value = 'Región de Ñuble'
$ws.Range("P9").Value = 83
$ws.Range("Q9").Value = 6

$ws.Range("D10").Value = 44525
$ws.Range("I10").Value = 'Primera'
$ws.Range("J10").Value = 200
$ws.Range("K10").Value = 600
$ws.Range("L10").Value = 700
$ws.Range("M10").Value = 650
$ws.Range("N10").Value = '$/paquete 6 unidades'
$ws.Range("O10").Value = 'Región de Ñuble'
$ws.Range("P10").Value = 108
$ws.Range("Q10").Value = 6

$ws.Range("D11").Value = 44525
$ws.Range("I11").Value = 'Segunda'
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 500
$ws.Range("L11").Value = 500
$ws.Range("M11").Value = 500
$ws.Range("N11").Value = '$/paquete 6 unidades'
$ws.Range("O11").Value = 'Región de Ñuble'
$ws.Range("P11").Value = 83
$ws.Range("Q11").Value = 6

$ws.Range("D12").Value = 44637
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 110
$ws.Range("K12").Value = 6500
$ws.Range("L12").Value = 7000
$ws.Range("M12").Value = 6773
$ws.Range("N12").Value = '$/paquete 36 unidades'
$ws.Range("O12").Value = 'Región Metropolitana'
$ws.Range("P12").Value = 188
$ws.Range("Q12").Value = 36

$ws.Range("D13").Value = 44328
$ws.Range("I13").Value = 'Primera'
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 600
$ws.Range("L13").Value = 700
$ws.Range("M13").Value = 650
$ws.Range("N13").Value = '$/paquete 6 unidades'
$ws.Range("O13").Value = 'Región de Ñuble'
$ws.Range("P13").Value = 108
$ws.Range("Q13").Value = 6

$ws.Range("D14").Value = 44328
$ws.Range("I14").Value = 'Segunda'
$ws.Range("J14").Value = 50
$ws.Range("K14").Value = 500
$ws.Range("L14").Value = 500
$ws.Range("M14").Value = 500
$ws.Range("N14").Value = '$/paquete 6 unidades'
$ws.Range("O14").Value = 'Región de Ñuble'
$ws.Range("P14").Value = 83
$ws.Range("Q14").Value = 6

$ws.Range("D15").Value = 44643
$ws.Range("I15").Value = 'Primera'
$ws.Range("J15").Value = 180
$ws.Range("K15").Value = 6500
$ws.Range("L15").Value = 7000
$ws.Range("M15").Value = 6778
$ws.Range("N15").Value = '$/paquete 36 unidades'
$ws.Range("O15").Value = 'Región Metropolitana'
$ws.Range("P15").Value = 188
$ws.Range("Q15").Value = 36

$ws.Range("D16").Value = 44188
$ws.Range("I16").Value = 'Primera'
$ws.Range("J16").Value = 200
$ws.Range("K16").Value = 600
$ws.Range("L16").Value = 700
$ws.Range("M16").Value = 650
$ws.Range("N16").Value = '$/paquete 6 unidades'
$ws.Range("O16").Value = 'Región de Ñuble'
$ws.Range("P16").Value = 108
$ws.Range("Q16").Value = 6

$ws.Range("D17").Value = 44188
$ws.Range("I17").Value = 'Segunda'
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 500
$ws.Range("L17").Value = 500
$ws.Range("M17").Value = 500
$ws.Range("N17").Value = '$/paquete 6 unidades'
$ws.Range("O17").Value = 'Región de Ñuble'
$ws.Range("P17").Value = 83
$ws.Range("Q17").Value = 6

$ws.Range("D18").Value = 44554
$ws.Range("I18").Value = 'Primera'
$ws.Range("J18").Value = 200
$ws.Range("K18").Value = 600
$ws.Range("L18").Value = 700
$ws.Range("M18").Value = 650
$ws.Range("N18").Value = '$/paquete 6 unidades'
$ws.Range("O18").Value = 'Región de Ñuble'
$ws.Range("P18").Value = 108
$ws.Range("Q18").Value = 6

$ws.Range("D19").Value = 44554
$ws.Range("I19").Value = 'Segunda'
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 500
$ws.Range("L19").Value = 500
$ws.Range("M19").Value = 500
$ws.Range("N19").Value = '$/paquete 6 unidades'
$ws.Range("O19").Value = 'Región de Ñuble'
$ws.Range("P19").Value = 83
$ws.Range("Q19").Value = 6

$ws.Range("D20").Value = 44649
$ws.Range("I20").Value = 'Primera'
$ws.Range("J20").Value = 220
$ws.Range("K20").Value = 8000
$ws.Range("L20").Value = 8500
$ws.Range("M20").Value = 8227
$ws.Range("N20").Value = '$/paquete 36 unidades'
$ws.Range("O20").Value = 'Región Metropolitana'
$ws.Range("P20").Value = 229
$ws.Range("Q20").Value = 36

$ws.Range("D21").Value = 44491
$ws.Range("I21").Value = 'Primera'
$ws.Range("J21").Value = 200
$ws.Range("K21").Value = 600
$ws.Range("L21").Value = 700
$ws.Range("M21").Value = 650
$ws.Range("N21").Value = '$/paquete 6 unidades'
$ws.Range("O21").Value = 'Región Metropolitana'
$ws.Range("P21").Value = 108
$ws.Range("Q21").Value = 6

$ws.Range("D22").Value = 44491
$ws.Range("I22").Value = 'Segunda'
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = 500
$ws.Range("N22").Value = '$/paquete 6 unidades'
$ws.Range("O22").Value = 'Región Metropolitana'
$ws.Range("P22").Value = 83
$ws.Range("Q22").Value = 6

$ws.Range("D23").Value = 44595
$ws.Range("I23").Value = 'Primera'
$ws.Range("J23").Value = 200
$ws.Range("K23").Value = 600
$ws.Range("L23").Value = 700
$ws.Range("M23").Value = 650
$ws.Range("N23").Value = '$/paquete 6 unidades'
$ws.Range("O23").Value = 'Región Metropolitana'
$ws.Range("P23").Value = 108
$ws.Range("Q23").Value = 6

$ws.Range("D24").Value = 44293
$ws.Range("I24").Value = 'Primera'
$ws.Range("J24").Value = 100
$ws.Range("K24").Value = 600
$ws.Range("L24").Value = 700
$ws.Range("M24").Value = 650
$ws.Range("N24").Value = '$/paquete 6 unidades'
$ws.Range("O24").Value = 'Región de Ñuble'
$ws.Range("P24").Value = 108
$ws.Range("Q24").Value = 6

$ws.Range("D25").Value = 44293
$ws.Range("I25").Value = 'Segunda'
$ws.Range("J25").Value = 50
$ws.Range("K25").Value = 500
$ws.Range("L25").Value = 500
$ws.Range("M25").Value = 500
$ws.Range("N25").Value = '$/paquete 6 unidades'
$ws.Range("O25").Value = 'Región de Ñuble'
$ws.Range("P25").Value = 83
$ws.Range("Q25").Value = 6

$ws.Range("D26").Value = 44644
$ws.Range("I26").Value = 'Primera'
$ws.Range("J26").Value = 160
$ws.Range("K26").Value = 6500
$ws.Range("L26").Value = 7000
$ws.Range("M26").Value = 6750
$ws.Range("N26").Value = '$/paquete 36 unidades'
$ws.Range("O26").Value = 'Región Metropolitana'
$ws.Range("P26").Value = 188
$ws.Range("Q26").Value = 36

$ws.Range("D27").Value = 44616
$ws.Range("I27").Value = 'Primera'
$ws.Range("J27").Value = 200
$ws.Range("K27").Value = 600
$ws.Range("L27").Value = 700
$ws.Range("M27").Value = 650
$ws.Range("N27").Value = '$/paquete 6 unidades'
$ws.Range("O27").Value = 'Región de Ñuble'
$ws.Range("P27").Value = 108
$ws.Range("Q27").Value = 6

$ws.Range("D28").Value = 44616
$ws.Range("I28").Value = 'Segunda'
$ws.Range("J28").Value = 100
$ws.Range("K28").Value = 500
$ws.Range("L28").Value = 500
$ws.Range("M28").Value = 500
$ws.Range("N28").Value = '$/paquete 6 unidades'
$ws.Range("O28").Value = 'Región de Ñuble'
$ws.Range("P28").Value = 83
$ws.Range("Q28").Value = 6

$ws.Range("D29").Value = 44335
$ws.Range("I29").Value = 'Primera'
$ws.Range("J29").Value = 150
$ws.Range("K29").Value = 600
$ws.Range("L29").Value = 700
$ws.Range("M29").Value = 633
$ws.Range("N29").Value = '$/paquete 6 unidades'
$ws.Range("O29").Value = 'Región de Ñuble'
$ws.Range("P29").Value = 106
$ws.Range("Q29").Value = 6

$ws.Range("D30").Value = 44335
$ws.Range("I30").Value = 'Segunda'
$ws.Range("J30").Value = 50
$ws.Range("K30").Value = 500
$ws.Range("L30").Value = 500
$ws.Range("M30").Value = 500
$ws.Range("N30").Value = '$/paquete 6 unidades'
$ws.Range("O30").Value = 'Región de Ñuble'
$ws.Range("P30").Value = 83
$ws.Range("Q30").Value = 6

$ws.Range("D31").Value = 44358
$ws.Range("I31").Value = 'Primera'
$ws.Range("J31").Value = 200
$ws.Range("K31").Value = 600
$ws.Range("L31").Value = 700
$ws.Range("M31").Value = 650
$ws.Range("N31").Value = '$/paquete 6 unidades'
$ws.Range("O31").Value = 'Región de Ñuble'
$ws.Range("P31").Value = 108
$ws.Range("Q31").Value = 6

# New row 32
$ws.Range("A32").Value = 11
$ws.Range("B32").Value = 'Vega Monumental Concepción'
$ws.Range("C32").Value = 'Bíobío'
$ws.Range("D32").Value = 44358
$ws.Range("D32").NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Range("E32").Value = 8
$ws.Range("F32").Value = 100112037
$ws.Range("G32").Value = 'Cebollín'
$ws.Range("H32").Value = 'Sin especificar'
$ws.Range("I32").Value = 'Segunda'
$ws.Range("J32").Value = 100
$ws.Range("K32").Value = 500
$ws.Range("L32").Value = 500
$ws.Range("M32").Value = 500
$ws.Range("N32").Value = '$/paquete 6 unidades'
$ws.Range("O32").Value = 'Región de Ñuble'
$ws.Range("P32").Value = 83
$ws.Range("Q32").Value = 6
$ws.Range("R32").Value = 'Hortaliza'
